# Insert one new data row above current row 85 (Thursday 2021-09-08 market entry),
# pushing the existing rows 85-138 down to 86-139. This matches the dimension
# change A1:R138 -> A1:R139 and the row-by-row shift visible in the diff
# (new row 86 == old row 85, new row 87 == old row 86, ..., new row 139 ==
# old row 138), with the brand-new record landing in row 85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 85:138 down to 86:139, leaving a blank row 85 (inherits the
# date-format style from the row above, same as Excel's native Insert).
$ws.Rows.Item(85).Insert()

# Populate the new row 85 with the new "Pepino ensalada" market record.
$ws.Cells.Item(85, 1).Value = 7
$ws.Cells.Item(85, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(85, 3).Value = "Ñuble"
$ws.Cells.Item(85, 4).Value = 44447
$ws.Cells.Item(85, 5).Value = 16
$ws.Cells.Item(85, 6).Value = 100112043
$ws.Cells.Item(85, 7).Value = "Pepino ensalada"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 160
$ws.Cells.Item(85, 11).Value = 17000
$ws.Cells.Item(85, 12).Value = 18000
$ws.Cells.Item(85, 13).Value = 17500
$ws.Cells.Item(85, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 292
$ws.Cells.Item(85, 17).Value = 60
$ws.Cells.Item(85, 18).Value = "Hortaliza"
